# Add season-record columns (Wins, Losses, Ties) to the roster/statistics
# table, matching the header style already used by the other column
# headers in row 1, then fill every player row with the team's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of an existing header cell (bold, centered,
# bordered) onto the new header cells so they match the rest of row 1
# without introducing new style entries.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows ---------------------------------------------------------
# Every player on the roster shares the team's overall season record.
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 90  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 72  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
